$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so numeric-looking strings
# (e.g. "182.70", "1.00", "66.413.11") are not auto-coerced to numbers,
# matching the source workbook where these are stored as text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '66.413.11'
$ws.Range('E2').Value = '  -4.22%  '
$ws.Range('D3').Value = '3.358.64'
$ws.Range('E3').Value = '  -4.77%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '559.95'
$ws.Range('E5').Value = '  -3.87%  '
$ws.Range('D6').Value = '182.70'
$ws.Range('E6').Value = '  -6.83%  '
$ws.Range('D7').Value = '0.598'
$ws.Range('E7').Value = '  -2.11%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '3.349.91'
$ws.Range('E9').Value = '  -4.72%  '
$ws.Range('D10').Value = '0.186'
$ws.Range('E10').Value = '  -8.28%  '
$ws.Range('D11').Value = '0.592'
$ws.Range('E11').Value = '  -5.08%  '
$ws.Range('D12').Value = '47.87'
$ws.Range('E12').Value = '  -7.44%  '
$ws.Range('D13').Value = '0.0000268'
$ws.Range('E13').Value = '  -5.92%  '
$ws.Range('D14').Value = '8.70'
$ws.Range('E14').Value = '  -5.92%  '
$ws.Range('D15').Value = '3.895.18'
$ws.Range('E15').Value = '  -4.88%  '
$ws.Range('D16').Value = '607.44'
$ws.Range('E16').Value = '  -10.10%  '
$ws.Range('D17').Value = '18.18'
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').Value = '66.334.72'
$ws.Range('E18').Value = '  -4.55%  '
$ws.Range('D19').Value = '3.360.84'
$ws.Range('E19').Value = '  -5.13%  '
$ws.Range('D20').Value = '0.117'
$ws.Range('E20').Value = '  -3.26%  '
$ws.Range('D21').Value = '11.48'
$ws.Range('E21').Value = '  -7.57%  '
$ws.Range('D22').Value = '0.912'
$ws.Range('E22').Value = '  -5.63%  '
$ws.Range('D23').Value = '16.83'
$ws.Range('E23').Value = '  -8.24%  '
$ws.Range('E24').Value = '  -1.70%  '
$ws.Range('D25').Value = '99.95'
$ws.Range('E25').Value = '  -5.79%  '
$ws.Range('D26').Value = '4.07'
$ws.Range('E26').Value = '  -6.80%  '
$ws.Range('D27').Value = '6.01'
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('E28').Value = '  -7.63%  '
$ws.Range('D29').Value = '9.38'
$ws.Range('E29').Value = '  -7.89%  '
$ws.Range('D30').Value = '8.79'
$ws.Range('E30').Value = '  -9.25%  '
$ws.Range('D31').Value = '30.63'
$ws.Range('E31').Value = '  -7.90%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '6.30'
$ws.Range('E32').Value = '  -8.04%  '
$ws.Range('B33').Value = 'dogwifhat'
$ws.Range('C33').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D33').Value = '3.83'
$ws.Range('E33').Value = '  -12.63%  '
$ws.Range('D34').Value = '11.11'
$ws.Range('E34').Value = '  -6.18%  '
$ws.Range('D35').Value = '548.22'
$ws.Range('E35').Value = '  +9.77%  '
$ws.Range('D36').Value = '0.105'
$ws.Range('E36').Value = '  -5.03%  '
$ws.Range('D37').Value = '3.827.72'
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').Value = '58.08'
$ws.Range('E38').Value = '  -6.22%  '
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('E40').Value = '  -7.39%  '
$ws.Range('D41').Value = '0.0₃0718'
$ws.Range('E41').Value = '  -11.03%  '
$ws.Range('D42').Value = '3.46'
$ws.Range('E42').Value = '  +25.20%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '0.127'
$ws.Range('E43').Value = '  -5.32%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '2.66'
$ws.Range('E44').Value = '  -9.02%  '
$ws.Range('D45').Value = '0.347'
$ws.Range('E45').Value = '  -6.70%  '
$ws.Range('D46').Value = '32.21'
$ws.Range('E46').Value = '  -6.83%  '
$ws.Range('D47').Value = '0.0415'
$ws.Range('E47').Value = '  -9.48%  '
$ws.Range('D48').Value = '3.17'
$ws.Range('E48').Value = '  -6.55%  '
$ws.Range('D49').Value = '2.65'
$ws.Range('E49').Value = '  -8.54%  '
$ws.Range('E50').Value = '  -4.46%  '
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.15%  '

# Restore default (unstyled) cell style on the data range so only the
# cell contents change, keeping styling identical to the original file.
$ws.Range("D2:E51").Style = "Normal"
